$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

# Copy formatting from an existing header cell (A1) to the new header cells
$ws.Range("A1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Update existing numeric values in row 2
$ws.Range("B2").Value = 0.09785757333711262
$ws.Range("C2").Value = 0.9994170903566206
$ws.Range("D2").Value = 0.2195837108323081

# Update F2 text (multi-line)
$ws.Range("F2").Value = "Pipeline(steps=[('model',`n                 RandomForestRegressor(max_depth=5, n_estimators=150))])"

# New data cells
$ws.Range("G2").Value = 0.1256850772835605
$ws.Range("H2").Value = 0.99
